$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "256.37"
Set-TextValue $ws "E2" "-0.13%"
Set-TextValue $ws "G2" "11"

Set-TextValue $ws "D3" "26.83"
Set-TextValue $ws "E3" "0.08%"
Set-TextValue $ws "G3" "11"

Set-TextValue $ws "D4" "4.730"
Set-TextValue $ws "E4" "-0.27%"
Set-TextValue $ws "G4" "11"

Set-TextValue $ws "D5" "0.05944"
Set-TextValue $ws "E5" "0.29%"
Set-TextValue $ws "G5" "11"

Set-TextValue $ws "D6" "6.621"
Set-TextValue $ws "E6" "-0.61%"
Set-TextValue $ws "G6" "11"

Set-TextValue $ws "D7" "0.8499"
Set-TextValue $ws "E7" "-1.94%"
Set-TextValue $ws "G7" "11"

Set-TextValue $ws "D8" "0.9228"
Set-TextValue $ws "E8" "-2.28%"
Set-TextValue $ws "G8" "11"

Set-TextValue $ws "E9" "-1.58%"
Set-TextValue $ws "G9" "11"

Set-TextValue $ws "D10" "0.04215"
Set-TextValue $ws "E10" "10.87%"
Set-TextValue $ws "G10" "11"

Set-TextValue $ws "E11" "-1.46%"
Set-TextValue $ws "G11" "11"

Set-TextValue $ws "D12" "0.03053"
Set-TextValue $ws "E12" "-3.49%"
Set-TextValue $ws "G12" "11"

Set-TextValue $ws "D13" "0.09099"
Set-TextValue $ws "E13" "-1.58%"
Set-TextValue $ws "G13" "11"

Set-TextValue $ws "D14" "0.001537"
Set-TextValue $ws "E14" "-0.63%"
Set-TextValue $ws "G14" "11"

Set-TextValue $ws "D15" "0.0006032"
Set-TextValue $ws "E15" "-0.25%"
Set-TextValue $ws "G15" "11"

Set-TextValue $ws "D16" "0.006075"
Set-TextValue $ws "E16" "-0.07%"
Set-TextValue $ws "G16" "11"

Set-TextValue $ws "D17" "3.470"
Set-TextValue $ws "E17" "-0.84%"
Set-TextValue $ws "G17" "11"

Set-TextValue $ws "D18" "3.164"
Set-TextValue $ws "E18" "-1.16%"
Set-TextValue $ws "G18" "11"

Set-TextValue $ws "E19" "-0.91%"
Set-TextValue $ws "G19" "11"

Set-TextValue $ws "D20" "0.3029"
Set-TextValue $ws "E20" "-3.11%"
Set-TextValue $ws "G20" "11"

Set-TextValue $ws "G21" "11"

Set-TextValue $ws "D22" "3.939"
Set-TextValue $ws "E22" "3.56%"
Set-TextValue $ws "G22" "11"

Set-TextValue $ws "D23" "0.04263"
Set-TextValue $ws "E23" "0.99%"
Set-TextValue $ws "G23" "11"

Set-TextValue $ws "D24" "0.001222"
Set-TextValue $ws "E24" "-0.22%"
Set-TextValue $ws "G24" "11"

Set-TextValue $ws "E25" "-15.56%"
Set-TextValue $ws "G25" "11"

Set-TextValue $ws "D26" "0.0001201"
Set-TextValue $ws "E26" "0.16%"
Set-TextValue $ws "G26" "11"

Set-TextValue $ws "D27" "0.0001524"
Set-TextValue $ws "E27" "1.99%"
Set-TextValue $ws "G27" "11"

Set-TextValue $ws "G28" "11"

Set-TextValue $ws "G29" "11"

Set-TextValue $ws "G30" "11"

Set-TextValue $ws "G31" "11"

Set-TextValue $ws "G32" "11"

Set-TextValue $ws "G33" "11"

Set-TextValue $ws "G34" "11"

Set-TextValue $ws "G35" "11"

Set-TextValue $ws "G36" "11"

Set-TextValue $ws "G37" "11"

Set-TextValue $ws "G38" "11"

Set-TextValue $ws "G39" "11"

Set-TextValue $ws "D40" "0.03776"
Set-TextValue $ws "E40" "-1.34%"
Set-TextValue $ws "G40" "11"

Set-TextValue $ws "D41" "0.006287"
Set-TextValue $ws "E41" "0.81%"
Set-TextValue $ws "G41" "11"

Set-TextValue $ws "D42" "0.1097"
Set-TextValue $ws "E42" "-0.35%"
Set-TextValue $ws "G42" "11"

Set-TextValue $ws "D43" "0.002443"
Set-TextValue $ws "E43" "11.10%"
Set-TextValue $ws "G43" "11"

Set-TextValue $ws "D44" "0.01383"
Set-TextValue $ws "E44" "23.97%"
Set-TextValue $ws "G44" "11"

Set-TextValue $ws "D45" "0.00005324"
Set-TextValue $ws "E45" "-3.15%"
Set-TextValue $ws "G45" "11"

Set-TextValue $ws "E46" "0.07%"
Set-TextValue $ws "G46" "11"

Set-TextValue $ws "D47" "0.04401"
Set-TextValue $ws "E47" "-50.25%"
Set-TextValue $ws "G47" "11"

Set-TextValue $ws "E48" "9,837.78%"
Set-TextValue $ws "G48" "11"

Set-TextValue $ws "E49" "0.07%"
Set-TextValue $ws "G49" "11"

Set-TextValue $ws "D50" "0.0002000"
Set-TextValue $ws "E50" "0.07%"
Set-TextValue $ws "G50" "11"

Set-TextValue $ws "G51" "11"

Write-Host "Done"